$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Sundar Pichai 1.jpg", 0),
    @("Mark 6.png", 300),
    @("Mark3.png", 542),
    @("Sundar Pichai 3.png", 1112),
    @("satya nadella 2.jpeg", 1230),
    @("mark zuckerberg.jpeg", 1364),
    @("Mark 4.jpg", 1446),
    @("Mark1.jpeg", 2108),
    @("Mark2.jpg", 3494),
    @("satya nadella 2.png", 4450),
    @("Sundar Pichai 2.jpg", 115087)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
